$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-17 Thursday", "2024-10-18 Friday"),
    @("32×38=1216", "60×14=840"),
    @("18×85=1530", "26×13=338"),
    @("69×67=4623", "35×25=875"),
    @("18×31=558", "32×41=1312"),
    @("13×94=1222", "35×50=1750"),
    @("71×39=2769", "18×73=1314"),
    @("34×23=782", "79×64=5056"),
    @("83×71=5893", "28×71=1988"),
    @("26×81=2106", "14×55=770"),
    @("89×86=7654", "18×60=1080"),
    @("46×88=4048", "56×88=4928"),
    @("58×31=1798", "61×95=5795"),
    @("15×18=270", "89×92=8188"),
    @("78×89=6942", "74×34=2516"),
    @("63×47=2961", "73×59=4307"),
    @("53×91=4823", "88×80=7040"),
    @("56×85=4760", "20×32=640"),
    @("13×50=650", "23×96=2208"),
    @("42×86=3612", "37×71=2627"),
    @("59×52=3068", "83×36=2988"),
    @("15×75=1125", "73×51=3723"),
    @("35×12=420", "44×37=1628"),
    @("14×56=784", "87×82=7134"),
    @("42×52=2184", "12×76=912"),
    @("22×88=1936", "85×85=7225")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done"
